$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("BACKEND")
$ws2 = $wb.Worksheets.Item("FRONTEND")

# --- New shared strings / new rows (ordered to reproduce shared string table order) ---
$ws2.Range("A16").Value = "faire un readme pour starter le projet"
$ws1.Range("A14").Value = "faire les commentaires"
$ws1.Range("A15").Value = "vérifier les messages d'erreur"
$ws2.Range("A19").Value = "faire le tour du site pour trouver des erreurs"
$ws1.Range("A16").Value = "faire le fichier txt qui explique ce que chaque test fait"
$ws2.Range("A17").Value = "faire les commentaires"
$ws2.Range("A18").Value = "vérifier les messages d'erreur"

# --- Mark additional FRONTEND rows as done ("x" in column B) ---
$ws2.Range("B6").Value = "x"
$ws2.Range("B7").Value = "x"
$ws2.Range("B8").Value = "x"
$ws2.Range("B10").Value = "x"
$ws2.Range("B13").Value = "x"
$ws2.Range("B15").Value = "x"

# --- Selections / view state ---
$ws1.Range("A16").Select()
$ws2.Range("B13").Select()
$ws2.Activate()
